# Atualização de bases das ligas, do dia: 12-06-2024 às 23:38
# Swap the match-data (columns B:AD) between each pair of rows below.
# Column A (the running index) is left untouched in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(26, 27),
    @(155, 156),
    @(159, 160),
    @(183, 184),
    @(185, 186),
    @(190, 191),
    @(276, 277),
    @(313, 315)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}
